# Weekly price update: a new week's observation is inserted as a new row
# right before the existing row 62, pushing every subsequent row (old 62..197)
# down by one (new rows 63..198). The new row carries the latest date
# (44544) together with its own Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 62..197 down to 63..198, leaving a blank row 62 to populate.
$ws.Rows("62:62").Insert()

# Populate the newly inserted row 62 with this week's record.
$ws.Range("A62").Value = 9
$ws.Range("B62").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C62").Value = "Metropolitana"
$ws.Range("D62").Value = 44544
$ws.Range("E62").Value = 13
$ws.Range("F62").Value = 300000001
$ws.Range("G62").Value = "Rabanito"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 7900
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = 2747
$ws.Range("N62").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O62").Value = "Provincia de Chacabuco"
$ws.Range("P62").Value = 27
$ws.Range("Q62").Value = 100
$ws.Range("R62").Value = "Hortaliza"
